# "filled out input_variables where blank and col_id to id"
#
# Column F = "input_variables".
#  - F2 currently holds "col_id" -> change it to "id".
#  - Every other row in F that is currently blank gets filled with
#    "impossible" (matching the neighbouring Mlstr_harmo::rule_category /
#    algorithm / status columns, which already say "impossible" for those
#    rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# col_id -> id
$ws.Range("F2").Value = "id"

# Rows whose "input_variables" (column F) cell is blank in the original
# workbook; fill them in with "impossible".
$rowsToFill = @(5, 7, 8, 9, 10, 11, 12, 13, 14, 18, 19, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61)

foreach ($r in $rowsToFill) {
    $ws.Cells.Item($r, 6).Value = "impossible"
}

# Reflect the new active selection recorded in the saved view.
$ws.Activate()
$ws.Range("E2").Select()
